$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 365

# Header for the new column G
$ws.Cells.Item(1, 7).Value = "ema"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Pull the Close column (E2:E365) in one shot, compute an exponential
# moving average (span = 7, i.e. alpha = 2/(span+1) = 0.25, not
# bias-adjusted - same as pandas' ewm(span=7, adjust=False).mean()),
# then push the results back out in one shot.
$closeRange = $ws.Range("E2:E$lastRow")
$closeVals = $closeRange.Value2

$alpha = 2.0 / (7.0 + 1.0)

$n = $lastRow - 1
$emaVals = New-Object 'object[,]' $n, 1
$prev = $closeVals[1, 1]
$emaVals[0, 0] = $prev
for ($i = 2; $i -le $n; $i++) {
    $cur = $closeVals[$i, 1]
    $prev = ($alpha * $cur) + ((1 - $alpha) * $prev)
    $emaVals[$i - 1, 0] = $prev
}

$ws.Range("G2:G$lastRow").Value = $emaVals
